# Add data for 2022-05-24: roll the "through" date from 05-15 to 05-16
# and bump the May / Total figures in the "2022" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-05-16"

# Update the header label for the 2022 column (column I, row 1).
$ws.Range("I1").Value = "2022 (through 05-16)"

# Update the May figure (row 6) and the Total figure (row 14) for 2022.
$ws.Range("I6").Value = 59
$ws.Range("I14").Value = 611
